$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6, column A currently uses the "date only" format (style 3); the
# daily-update run re-stamps it to the "date + time" format (style 2),
# matching rows 2-5, since it is no longer the most-recent row.
$ws.Range("A6").NumberFormat = $ws.Range("A2").NumberFormat

# Append today's row (row 7) with the new day's tallies.
$ws.Range("A7").Value = 45747
$ws.Range("B7").Value = 22
$ws.Range("C7").Value = 21
$ws.Range("D7").Value = 24

# The newest row keeps the "date only" format that row 6 used to have.
$ws.Range("A7").NumberFormat = "YYYY-MM-DD"
